# Quarterly income-statement update: the reporting window rolls forward by one
# quarter (fiscal Q4 1401, ending 1401/12) - drop the oldest quarter column (was D)
# and append the newest quarter as the new last column (M). All header labels,
# publish-date labels, and financial figures shift left by one column; the new
# quarter's freshly-reported figures are written into column M.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Overview")

# --- Row 8: quarter-period headers (D8:M8) ---
$ws.Cells.Item(8, 4).Value = "فصل سوم منتهی به 1399/09"
$ws.Cells.Item(8, 5).Value = "فصل چهارم منتهی به 1399/12"
$ws.Cells.Item(8, 6).Value = "فصل اول منتهی به 1400/03"
$ws.Cells.Item(8, 7).Value = "فصل دوم منتهی به 1400/06"
$ws.Cells.Item(8, 8).Value = "فصل سوم منتهی به 1400/09"
$ws.Cells.Item(8, 9).Value = "فصل چهارم منتهی به 1400/12"
$ws.Cells.Item(8, 10).Value = "فصل اول منتهی به 1401/03"
$ws.Cells.Item(8, 11).Value = "فصل دوم منتهی به 1401/06"
$ws.Cells.Item(8, 12).Value = "فصل سوم منتهی به 1401/09"
$ws.Cells.Item(8, 13).Value = "فصل چهارم منتهی به 1401/12"

# --- Row 9: publish-date headers (D9:M9) ---
$ws.Cells.Item(9, 4).Value = "1400-11-02 (2)"
$ws.Cells.Item(9, 5).Value = "1401-02-25 (12)"
$ws.Cells.Item(9, 6).Value = "1401-05-01 (3)"
$ws.Cells.Item(9, 7).Value = "1401-08-30 (4)"
$ws.Cells.Item(9, 8).Value = "1401-11-19 (3)"
$ws.Cells.Item(9, 9).Value = "1402-02-27 (12)"
# Leading apostrophe forces text (this literal would otherwise auto-parse as a date)
$ws.Cells.Item(9, 10).Value = "'1401-05-01"
$ws.Cells.Item(9, 11).Value = "1401-08-30 (2)"
$ws.Cells.Item(9, 12).Value = "1401-11-19 (2)"
$ws.Cells.Item(9, 13).Value = "1402-02-27 (3)"

# --- Data rows 11-27: shift values left one quarter, append new quarter value ---
# Row 11: Sales (فروش)
$ws.Cells.Item(11, 4).Value = 2190426
$ws.Cells.Item(11, 5).Value = 1987544
$ws.Cells.Item(11, 6).Value = 2671724
$ws.Cells.Item(11, 7).Value = 2575491
$ws.Cells.Item(11, 8).Value = 3090529
$ws.Cells.Item(11, 9).Value = 3013415
$ws.Cells.Item(11, 10).Value = 3415613
$ws.Cells.Item(11, 11).Value = 4721905
$ws.Cells.Item(11, 12).Value = 6328926
$ws.Cells.Item(11, 13).Value = 6690899

# Row 12: Cost of goods sold
$ws.Cells.Item(12, 4).Value = -830393
$ws.Cells.Item(12, 5).Value = -792083
$ws.Cells.Item(12, 6).Value = -1072057
$ws.Cells.Item(12, 7).Value = -1076741
$ws.Cells.Item(12, 8).Value = -1323640
$ws.Cells.Item(12, 9).Value = -1727472
$ws.Cells.Item(12, 10).Value = -1558599
$ws.Cells.Item(12, 11).Value = -1728304
$ws.Cells.Item(12, 12).Value = -2446043
$ws.Cells.Item(12, 13).Value = -3130883

# Row 13: Gross profit
$ws.Cells.Item(13, 4).Value = 1360033
$ws.Cells.Item(13, 5).Value = 1195461
$ws.Cells.Item(13, 6).Value = 1599667
$ws.Cells.Item(13, 7).Value = 1498750
$ws.Cells.Item(13, 8).Value = 1766889
$ws.Cells.Item(13, 9).Value = 1285943
$ws.Cells.Item(13, 10).Value = 1857014
$ws.Cells.Item(13, 11).Value = 2993601
$ws.Cells.Item(13, 12).Value = 3882883
$ws.Cells.Item(13, 13).Value = 3560016

# Row 14: G&A expenses
$ws.Cells.Item(14, 4).Value = -220301
$ws.Cells.Item(14, 5).Value = -432882
$ws.Cells.Item(14, 6).Value = -290181
$ws.Cells.Item(14, 7).Value = -190044
$ws.Cells.Item(14, 8).Value = -94076
$ws.Cells.Item(14, 9).Value = -158693
$ws.Cells.Item(14, 10).Value = -147803
$ws.Cells.Item(14, 11).Value = -187732
$ws.Cells.Item(14, 12).Value = -709440
$ws.Cells.Item(14, 13).Value = -761923

# Row 15: Impairment expense
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 0
$ws.Cells.Item(15, 7).Value = 0
$ws.Cells.Item(15, 8).Value = 0
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 10).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 0

# Row 16: Other operating income/expense, net
$ws.Cells.Item(16, 4).Value = 61880
$ws.Cells.Item(16, 5).Value = -25759
$ws.Cells.Item(16, 6).Value = -1761
$ws.Cells.Item(16, 7).Value = -46075
$ws.Cells.Item(16, 8).Value = 48276
$ws.Cells.Item(16, 9).Value = -27987
$ws.Cells.Item(16, 10).Value = 1810
$ws.Cells.Item(16, 11).Value = -2150
$ws.Cells.Item(16, 12).Value = 84653
$ws.Cells.Item(16, 13).Value = 295552

# Row 17: Operating profit
$ws.Cells.Item(17, 4).Value = 1201612
$ws.Cells.Item(17, 5).Value = 736820
$ws.Cells.Item(17, 6).Value = 1307725
$ws.Cells.Item(17, 7).Value = 1262631
$ws.Cells.Item(17, 8).Value = 1721089
$ws.Cells.Item(17, 9).Value = 1099263
$ws.Cells.Item(17, 10).Value = 1711021
$ws.Cells.Item(17, 11).Value = 2803719
$ws.Cells.Item(17, 12).Value = 3258096
$ws.Cells.Item(17, 13).Value = 3093645

# Row 18: Finance costs
$ws.Cells.Item(18, 4).Value = -12804
$ws.Cells.Item(18, 5).Value = -9412
$ws.Cells.Item(18, 6).Value = -13454
$ws.Cells.Item(18, 7).Value = -12095
$ws.Cells.Item(18, 8).Value = -12842
$ws.Cells.Item(18, 9).Value = -10177
$ws.Cells.Item(18, 10).Value = -14861
$ws.Cells.Item(18, 11).Value = -13210
$ws.Cells.Item(18, 12).Value = -14188
$ws.Cells.Item(18, 13).Value = -13201

# Row 19: Other non-operating income/expense, net
$ws.Cells.Item(19, 4).Value = 150230
$ws.Cells.Item(19, 5).Value = 20454
$ws.Cells.Item(19, 6).Value = 105926
$ws.Cells.Item(19, 7).Value = 65576
$ws.Cells.Item(19, 8).Value = 176427
$ws.Cells.Item(19, 9).Value = 106302
$ws.Cells.Item(19, 10).Value = 250879
$ws.Cells.Item(19, 11).Value = 246395
$ws.Cells.Item(19, 12).Value = 286213
$ws.Cells.Item(19, 13).Value = 466691

# Row 20: Profit before tax from continuing ops
$ws.Cells.Item(20, 4).Value = 1339038
$ws.Cells.Item(20, 5).Value = 747862
$ws.Cells.Item(20, 6).Value = 1400197
$ws.Cells.Item(20, 7).Value = 1316112
$ws.Cells.Item(20, 8).Value = 1884674
$ws.Cells.Item(20, 9).Value = 1195388
$ws.Cells.Item(20, 10).Value = 1947039
$ws.Cells.Item(20, 11).Value = 3036904
$ws.Cells.Item(20, 12).Value = 3530121
$ws.Cells.Item(20, 13).Value = 3547135

# Row 21: Tax
$ws.Cells.Item(21, 4).Value = -193764
$ws.Cells.Item(21, 5).Value = 21851
$ws.Cells.Item(21, 6).Value = -99622
$ws.Cells.Item(21, 7).Value = -362256
$ws.Cells.Item(21, 8).Value = -63772
$ws.Cells.Item(21, 9).Value = -164055
$ws.Cells.Item(21, 10).Value = -252155
$ws.Cells.Item(21, 11).Value = -669939
$ws.Cells.Item(21, 12).Value = -454234
$ws.Cells.Item(21, 13).Value = -390193

# Row 22: Net profit from continuing operations
$ws.Cells.Item(22, 4).Value = 1145274
$ws.Cells.Item(22, 5).Value = 769713
$ws.Cells.Item(22, 6).Value = 1300575
$ws.Cells.Item(22, 7).Value = 953856
$ws.Cells.Item(22, 8).Value = 1820902
$ws.Cells.Item(22, 9).Value = 1031333
$ws.Cells.Item(22, 10).Value = 1694884
$ws.Cells.Item(22, 11).Value = 2366965
$ws.Cells.Item(22, 12).Value = 3075887
$ws.Cells.Item(22, 13).Value = 3156942

# Row 23: Discontinued operations, net of tax
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(23, 9).Value = 0
$ws.Cells.Item(23, 10).Value = 0
$ws.Cells.Item(23, 11).Value = 0
$ws.Cells.Item(23, 12).Value = 0
$ws.Cells.Item(23, 13).Value = 0

# Row 24: Net profit
$ws.Cells.Item(24, 4).Value = 1145274
$ws.Cells.Item(24, 5).Value = 769713
$ws.Cells.Item(24, 6).Value = 1300575
$ws.Cells.Item(24, 7).Value = 953856
$ws.Cells.Item(24, 8).Value = 1820902
$ws.Cells.Item(24, 9).Value = 1031333
$ws.Cells.Item(24, 10).Value = 1694884
$ws.Cells.Item(24, 11).Value = 2366965
$ws.Cells.Item(24, 12).Value = 3075887
$ws.Cells.Item(24, 13).Value = 3156942

# Row 25: EPS after tax
$ws.Cells.Item(25, 4).Value = 1762
$ws.Cells.Item(25, 5).Value = 1184
$ws.Cells.Item(25, 6).Value = 2001
$ws.Cells.Item(25, 7).Value = 1467
$ws.Cells.Item(25, 8).Value = 2801
$ws.Cells.Item(25, 9).Value = 1587
$ws.Cells.Item(25, 10).Value = 2608
$ws.Cells.Item(25, 11).Value = 3641
$ws.Cells.Item(25, 12).Value = 4732
$ws.Cells.Item(25, 13).Value = 4857

# Row 26: Capital
$ws.Cells.Item(26, 4).Value = 650000
$ws.Cells.Item(26, 5).Value = 650000
$ws.Cells.Item(26, 6).Value = 650000
$ws.Cells.Item(26, 7).Value = 650000
$ws.Cells.Item(26, 8).Value = 650000
$ws.Cells.Item(26, 9).Value = 650000
$ws.Cells.Item(26, 10).Value = 650000
$ws.Cells.Item(26, 11).Value = 650000
$ws.Cells.Item(26, 12).Value = 650000
$ws.Cells.Item(26, 13).Value = 650000

# Row 27: EPS based on latest capital
$ws.Cells.Item(27, 4).Value = 1762
$ws.Cells.Item(27, 5).Value = 1184
$ws.Cells.Item(27, 6).Value = 2001
$ws.Cells.Item(27, 7).Value = 1467
$ws.Cells.Item(27, 8).Value = 2801
$ws.Cells.Item(27, 9).Value = 1587
$ws.Cells.Item(27, 10).Value = 2608
$ws.Cells.Item(27, 11).Value = 3641
$ws.Cells.Item(27, 12).Value = 4732
$ws.Cells.Item(27, 13).Value = 4857

# --- Column widths: shift left to match the new column layout (29/31 alternating) ---
$newWidths = @(28.17,30.17,28.17,28.17,28.17,30.17,28.17,28.17,28.17,30.17)
for ($i = 0; $i -lt 10; $i++) {
    $ws.Columns.Item(4 + $i).ColumnWidth = $newWidths[$i]
}
